$d = $word.ActiveDocument
$newText = "Perioadele campaniei din Pegasus: 8-17 octombrie, 7-16 noiembrie,"

# Collect every paragraph whose text still refers to the old 2018 Perseid
# campaign blurb before mutating anything (mutating while enumerating
# $d.Paragraphs directly would be unsafe).
$targets = New-Object System.Collections.ArrayList
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Perioadele campaniei*") {
        [void]$targets.Add($p)
    }
}

foreach ($p in $targets) {
    $start = $p.Range.Start
    $end = $p.Range.End
    # Exclude the trailing paragraph mark (last char of Paragraph.Range)
    # so only the runs inside the paragraph get removed.
    $body = $d.Range($start, $end - 1)
    $body.Delete()
    # Insert the translated sentence as a brand-new, plain run (no rPr).
    $ins = $d.Range($start, $start)
    $ins.InsertAfter($newText)
}
